$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab
$ws.Name = "Gamma2F"

# Add new row 16 data
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.013849313885671
$ws.Range("D16").Value = 0.9229541462827703
$ws.Range("E16").Value = 1.011873886281002
$ws.Range("F16").Value = 1.013849313885671
$ws.Range("G16").Value = 0.960652774054112
$ws.Range("H16").Value = 1.041031941877974
$ws.Range("I16").Value = 1.012512863396901
$ws.Range("J16").Value = 0.9229541462827703
$ws.Range("K16").Value = 0.9674140162818861
$ws.Range("L16").Value = 0.9906316650837783
$ws.Range("M16").Value = 0.9938124876297384

# Apply same style as A15 (border/bold/centered) to A16
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122) # xlPasteFormats
